$en0 = @"
After moving out of our old house, we moved into a larger flat.
"@
$en1 = @"
Before moving out, we need to get our new house ready for us.
"@
$en2 = @"
They'll be moving into their new house as soon as they've finished fitting out the living room.
"@
$en3 = @"
Moving to a new city will allow us to discover more about new cultures and traditions.
"@
$en4 = @"
By moving to this area, we hope to meet more friendly neighbours.
"@
$en5 = @"
We need to see more of each other, because I see a lot of advantages in that.
"@
$en6 = @"
The current situation has more disadvantages than advantages.
"@
$de0 = @"
Nachdem wir aus unserem alten Haus ausgezogen sind, sind wir in eine geräumigere Wohnung gezogen.
"@
$de1 = @"
Bevor wir umziehen, müssen wir unser neues Haus einrichten, damit es bereit ist, uns aufzunehmen.
"@
$de2 = @"
Sie werden in ihr neues Haus einziehen, sobald sie das Wohnzimmer eingerichtet haben.
"@
$de3 = @"
Der Umzug in eine neue Stadt ermöglicht es uns, mehr von den neuen Kulturen und Traditionen kennenzulernen.
"@
$de4 = @"
Durch den Umzug in diese Gegend hoffen wir, mehr nette Nachbarn zu treffen.
"@
$de5 = @"
Wir sollten uns mehr treffen, weil ich viele Vorteile darin sehe.
"@
$de6 = @"
Die derzeitige Situation hat mehr Nachteile als Vorteile.
"@
$nl0 = @"
Na het verlaten van ons oude huis zijn we verhuisd naar een grotere flat.
"@
$nl1 = @"
Voordat we verhuizen, moeten we ons nieuwe huis klaarmaken.
"@
$nl2 = @"
Ze verhuizen naar hun nieuwe huis zodra ze klaar zijn met het inrichten van de woonkamer.
"@
$nl3 = @"
Door te verhuizen naar een nieuwe stad kunnen we meer te weten komen over nieuwe culturen en tradities.
"@
$nl4 = @"
Door naar dit gebied te verhuizen, hopen we meer vriendelijke buren te ontmoeten.
"@
$nl5 = @"
We moeten elkaar meer gaan zien, want daar zie ik veel voordelen in.
"@
$nl6 = @"
De huidige situatie heeft meer nadelen dan voordelen.
"@
$pt0 = @"
Depois de sairmos da nossa antiga casa, mudámo-nos para um apartamento maior.
"@
$pt1 = @"
Antes de nos mudarmos, temos de preparar a nossa nova casa para nós.
"@
$pt2 = @"
Eles vão mudar-se para a nova casa assim que acabarem de equipar a sala de estar.
"@
$pt3 = @"
A mudança para uma nova cidade permitir-nos-á descobrir mais sobre novas culturas e tradições.
"@
$pt4 = @"
Ao mudarmo-nos para esta zona, esperamos encontrar mais vizinhos simpáticos.
"@
$pt5 = @"
Precisamos de nos ver mais uns aos outros, porque vejo muitas vantagens nisso.
"@
$pt6 = @"
A situação atual tem mais desvantagens do que vantagens.
"@
$fr0 = @"
Après avoir déménagé de notre ancienne maison, nous avons emménagé dans un appartement plus spacieux.
"@
$fr1 = @"
Avant de déménager, nous devons aménager notre nouvelle maison pour qu'elle soit prête à nous accueillir.
"@
$fr2 = @"
Ils vont emménager dans leur nouvelle maison dès qu'ils auront fini d'aménager le salon.
"@
$fr3 = @"
Déménager dans une nouvelle ville nous permettra de découvrir davantage de nouvelles cultures et traditions.
"@
$fr4 = @"
En emménageant dans ce quartier, nous espérons rencontrer davantage de voisins sympathiques.
"@
$fr5 = @"
Il faudrait qu'on se voie davantage car j'y vois beaucoup d'avantages.
"@
$fr6 = @"
La situation actuelle présente plus d’inconvénients que d’avantages.
"@
$es0 = @"
Después de mudarnos de nuestra antigua casa, nos mudamos a un piso más grande.
"@
$es1 = @"
Antes de mudarnos, tenemos que preparar nuestra nueva casa.
"@
$es2 = @"
Se mudarán a su nueva casa en cuanto terminen de acondicionar el salón.
"@
$es3 = @"
Mudarnos a una nueva ciudad nos permitirá descubrir nuevas culturas y tradiciones.
"@
$es4 = @"
Al mudarnos a esta zona, esperamos conocer a más vecinos amistosos.
"@
$es5 = @"
Tenemos que vernos más, porque le veo muchas ventajas.
"@
$es6 = @"
La situación actual tiene más desventajas que ventajas.
"@
$it0 = @"
Dopo aver lasciato la nostra vecchia casa, ci siamo trasferiti in un appartamento più grande.
"@
$it1 = @"
Prima di traslocare, dobbiamo preparare la nuova casa per noi.
"@
$it2 = @"
Si trasferiranno nella loro nuova casa non appena avranno finito di sistemare il soggiorno.
"@
$it3 = @"
Trasferirci in una nuova città ci permetterà di scoprire nuove culture e tradizioni.
"@
$it4 = @"
Trasferendoci in questa zona, speriamo di incontrare più vicini amichevoli.
"@
$it5 = @"
Dobbiamo vederci di più, perché vedo molti vantaggi in questo.
"@
$it6 = @"
La situazione attuale presenta più svantaggi che vantaggi.
"@
$note1 = @"
Les termes "déménager," "emménager," et "aménager" sont tous liés à des actions concernant des lieux de résidence ou des espaces, mais ils ont des significations distinctes.

#Déménager#
#Définition :# Déménager signifie quitter un logement pour s'installer dans un autre. Cela implique de transporter ses affaires d'un endroit à un autre.
#Exemple :#  "Nous allons déménager de notre appartement actuel pour une maison plus grande."

#Emménager#
#Définition :# Emménager signifie s'installer dans un nouveau logement. C'est l'action de prendre possession d'un nouvel espace de vie.
#Exemple :# "Nous avons emménagé dans notre nouvelle maison la semaine dernière."

#Aménager#
#Définition :# Aménager signifie organiser, arranger ou équiper un espace pour le rendre fonctionnel ou agréable à vivre. Cela peut inclure des travaux de rénovation, la décoration, ou l'installation de meubles.
#Exemple :# "Nous avons aménagé le grenier pour en faire une chambre d'amis."
"@
$note2 = @"
Les termes "davantage" et "d'avantage" sont souvent confondus, mais ils ont des significations et des usages distincts.

#Davantage#
#Définition :# "Davantage" est un adverbe qui signifie "plus" ou "en plus grande quantité." Il est utilisé pour indiquer une augmentation ou une quantité supplémentaire.
#Exemples :#
"Je voudrais davantage de temps pour terminer ce projet."
"Il a besoin de davantage d'informations pour prendre une décision."
"Elle aimerait davantage de soutien de la part de ses collègues."

#D'avantage#
#Définition :# "D'avantage" est une expression qui signifie "de bénéfice" ou "de profit." Elle est utilisée pour parler d'un avantage ou d'un bénéfice.
#Exemples :#
"Cette nouvelle technologie n'apporte pas d'avantage par rapport à l'ancienne."
"Il n'y a pas d'avantage à changer de fournisseur pour le moment."
"Cette solution présente plusieurs d'avantages par rapport aux autres options."

#Différence clé#
"Davantage" est utilisé pour indiquer une quantité supplémentaire ou une augmentation.
"D'avantage" est utilisé pour parler d'un bénéfice ou d'un avantage.
"@

$wb = $excel.ActiveWorkbook
$wsSentences = $wb.Worksheets.Item("Sentences")
$wsNotes = $wb.Worksheets.Item("Notes")

# --- Fix two Notes entries: trailing " :" -> "." before the blank line ---
$wsNotes.Range("C2").Value2 = $note1
$wsNotes.Range("C3").Value2 = $note2

# --- Append 7 new sentence rows (IDs 101-107) to the Sentences sheet ---
$ids = @(101, 102, 103, 104, 105, 106, 107)
$frArr = @($fr0, $fr1, $fr2, $fr3, $fr4, $fr5, $fr6)
$enArr = @($en0, $en1, $en2, $en3, $en4, $en5, $en6)
$deArr = @($de0, $de1, $de2, $de3, $de4, $de5, $de6)
$nlArr = @($nl0, $nl1, $nl2, $nl3, $nl4, $nl5, $nl6)
$ptArr = @($pt0, $pt1, $pt2, $pt3, $pt4, $pt5, $pt6)
$esArr = @($es0, $es1, $es2, $es3, $es4, $es5, $es6)
$itArr = @($it0, $it1, $it2, $it3, $it4, $it5, $it6)

$startRow = 102

# Column A (IDs) first
for ($i = 0; $i -lt 7; $i++) {
    $wsSentences.Cells.Item($startRow + $i, 1).Value2 = $ids[$i]
}
# Then column F (FR / source language) for all new rows
for ($i = 0; $i -lt 7; $i++) {
    $wsSentences.Cells.Item($startRow + $i, 6).Value2 = $frArr[$i]
}
# Then column B (EN)
for ($i = 0; $i -lt 7; $i++) {
    $wsSentences.Cells.Item($startRow + $i, 2).Value2 = $enArr[$i]
}
# Then column C (DE)
for ($i = 0; $i -lt 7; $i++) {
    $wsSentences.Cells.Item($startRow + $i, 3).Value2 = $deArr[$i]
}
# Then column D (NL)
for ($i = 0; $i -lt 7; $i++) {
    $wsSentences.Cells.Item($startRow + $i, 4).Value2 = $nlArr[$i]
}
# Then column E (PT)
for ($i = 0; $i -lt 7; $i++) {
    $wsSentences.Cells.Item($startRow + $i, 5).Value2 = $ptArr[$i]
}
# Then column G (ES)
for ($i = 0; $i -lt 7; $i++) {
    $wsSentences.Cells.Item($startRow + $i, 7).Value2 = $esArr[$i]
}
# Then column H (IT)
for ($i = 0; $i -lt 7; $i++) {
    $wsSentences.Cells.Item($startRow + $i, 8).Value2 = $itArr[$i]
}

# --- Notes sheet: selection resets to B3, no longer the active tab ---
$wsNotes.Activate() | Out-Null
$wsNotes.Range("B3").Select() | Out-Null

# --- View state: Sentences sheet becomes active, scrolled near the bottom, F108 selected ---
$wsSentences.Activate() | Out-Null
$wsSentences.Range("F108").Select() | Out-Null
